$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before T (Custom Field 1 / Custom field 2 shift right into U / V)
$ws.Columns("T:T").Insert()

# 2. Header for the newly inserted column
$ws.Range("T1").Value = "Form Tag"

# 3. Fill the new column's data rows with "Default", copying the (hyperlink-like) number
#    format/font that the source workbook used for this column (matches column S's style).
$ws.Range("S2").Copy()
$ws.Range("T2:T5").PasteSpecial(-4122)
$ws.Range("T2").Value = "Default"
$ws.Range("T3").Value = "Default"
$ws.Range("T4").Value = "Default"
$ws.Range("T5").Value = "Default"

# 4. Data correction on row 4 ("Verified *" column)
$ws.Range("N4").Value = "No"

# 5. The comment that documented the custom-field columns now belongs on the
#    (shifted) U1 cell instead of T1.
$comment = $ws.Range("T1").Comment
$commentText = $comment.Text()
$comment.Delete()
$ws.Range("U1").AddComment($commentText)

# 6. Restore the selection to reflect where the edit was made.
$ws.Range("T3:T5").Select()
